$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4
$ws.Range("B4").Value = -3105.68
$ws.Range("C4").Value = -146.3
$ws.Range("D4").Value = -3712.55
$ws.Range("E4").Value = -3054.75
$ws.Range("F4").Value = -1790.51
$ws.Range("G4").Value = -2040.45
$ws.Range("H4").Value = -4181.66
$ws.Range("I4").Value = -4518.96
$ws.Range("J4").Value = -2686.92
$ws.Range("K4").Value = -3115.19
$ws.Range("L4").Value = -4035.62
$ws.Range("M4").Value = -5904.47
$ws.Range("N4").Value = -3865.53
$ws.Range("O4").ClearContents()
$ws.Range("P4").Value = -749.02
$ws.Range("Q4").Value = -5790.97
$ws.Range("R4").Value = -1338.25
$ws.Range("S4").Value = -4269.16
$ws.Range("T4").Value = -2766.52
$ws.Range("U4").Value = -2425.15
$ws.Range("V4").Value = -4905.16
$ws.Range("W4").Value = -1749.82
$ws.Range("X4").Value = -1967.41
$ws.Range("Y4").Value = -2127.08
$ws.Range("Z4").Value = -70247.13

# Row 5
$ws.Range("B5").Value = -4826.34
$ws.Range("C5").Value = -2639.33
$ws.Range("D5").Value = -3018.89
$ws.Range("E5").Value = -2416.87
$ws.Range("F5").Value = -2223.36
$ws.Range("G5").Value = -1574.37
$ws.Range("H5").Value = -3896.58
$ws.Range("I5").Value = -5491.9
$ws.Range("J5").Value = -2642.31
$ws.Range("K5").Value = -2215.6
$ws.Range("L5").Value = -4605.03
$ws.Range("M5").Value = -3268.83
$ws.Range("N5").Value = -2129.77
$ws.Range("O5").Value = -4427.59
$ws.Range("P5").Value = -4021.79
$ws.Range("Q5").Value = -3562.04
$ws.Range("R5").Value = -5754.01
$ws.Range("S5").Value = -800.23
$ws.Range("T5").Value = -2138.08
$ws.Range("U5").Value = -3594.87
$ws.Range("V5").Value = -4849.86
$ws.Range("W5").Value = -2576.66
$ws.Range("X5").Value = -2616.12
$ws.Range("Y5").Value = -2475.07
$ws.Range("Z5").Value = -77765.5

# Row 6
$ws.Range("B6").Value = -1358.42
$ws.Range("C6").Value = -440.1
$ws.Range("D6").Value = -749.86
$ws.Range("E6").Value = -1555.01
$ws.Range("F6").Value = -180.29
$ws.Range("G6").Value = -211.18
$ws.Range("H6").Value = -1042.55
$ws.Range("I6").Value = -3922.16
$ws.Range("J6").Value = -2045.24
$ws.Range("K6").Value = -2789.62
$ws.Range("L6").Value = -60.54
$ws.Range("M6").Value = -1083.77
$ws.Range("N6").Value = -958.06
$ws.Range("O6").Value = -1319.77
$ws.Range("P6").Value = -803.88
$ws.Range("Q6").Value = -1156.78
$ws.Range("R6").Value = -1008.54
$ws.Range("S6").ClearContents()
$ws.Range("T6").Value = -2238.62
$ws.Range("U6").Value = -3379.73
$ws.Range("V6").Value = -3544.38
$ws.Range("W6").ClearContents()
$ws.Range("X6").Value = -190.59
$ws.Range("Y6").Value = -407.8
$ws.Range("Z6").Value = -30446.89

# Row 7
$ws.Range("B7").Value = -3983.31
$ws.Range("C7").Value = -3028.4
$ws.Range("D7").Value = -1627.97
$ws.Range("E7").Value = -6414.87
$ws.Range("F7").Value = -2669.13
$ws.Range("G7").Value = -2775.83
$ws.Range("H7").Value = -1135.31
$ws.Range("I7").Value = -1235.7
$ws.Range("J7").Value = -761.67
$ws.Range("K7").Value = -2445.23
$ws.Range("L7").Value = -872.38
$ws.Range("M7").Value = -4217.25
$ws.Range("N7").Value = -4137.41
$ws.Range("O7").Value = -2178.44
$ws.Range("P7").Value = -4309.05
$ws.Range("Q7").Value = -3071.57
$ws.Range("R7").Value = -3371.08
$ws.Range("S7").Value = -4151.35
$ws.Range("T7").Value = -3343.77
$ws.Range("U7").Value = -3395.33
$ws.Range("V7").Value = -3555.93
$ws.Range("W7").Value = -873.65
$ws.Range("X7").Value = -1850.57
$ws.Range("Y7").Value = -3554.36
$ws.Range("Z7").Value = -68959.56

# Row 8
$ws.Range("B8").ClearContents()
$ws.Range("C8").Value = -532.13
$ws.Range("D8").Value = -214.07
$ws.Range("E8").Value = -2470.26
$ws.Range("F8").Value = -4075.11
$ws.Range("G8").Value = -2366.48
$ws.Range("H8").Value = -1762.6
$ws.Range("I8").Value = -2865.75
$ws.Range("J8").Value = -705.02
$ws.Range("K8").Value = -1140.66
$ws.Range("L8").ClearContents()
$ws.Range("M8").ClearContents()
$ws.Range("N8").ClearContents()
$ws.Range("O8").Value = -2194
$ws.Range("P8").Value = -1654.62
$ws.Range("Q8").Value = -2074.53
$ws.Range("R8").Value = -765.97
$ws.Range("S8").Value = -4682.21
$ws.Range("T8").Value = -684.18
$ws.Range("U8").Value = -1584.88
$ws.Range("V8").Value = -1437.23
$ws.Range("W8").Value = -2007.79
$ws.Range("X8").Value = -48.04
$ws.Range("Y8").Value = -1486.74
$ws.Range("Z8").Value = -34752.27

# Row 9
$ws.Range("B9").Value = -5615.31
$ws.Range("C9").Value = -4276.73
$ws.Range("D9").Value = -2180.54
$ws.Range("E9").Value = -3693.87
$ws.Range("F9").Value = -3652.87
$ws.Range("G9").Value = -2784.62
$ws.Range("H9").Value = -4687.37
$ws.Range("I9").Value = -5621.78
$ws.Range("J9").Value = -3490
$ws.Range("K9").Value = -8178.43
$ws.Range("L9").Value = -4563.68
$ws.Range("M9").Value = -5282.9
$ws.Range("N9").Value = -6180.82
$ws.Range("O9").Value = -5229.73
$ws.Range("P9").Value = -7761.23
$ws.Range("Q9").Value = -4276.43
$ws.Range("R9").Value = -3012.77
$ws.Range("S9").Value = -5888.13
$ws.Range("T9").Value = -10311.62
$ws.Range("U9").Value = -5128.16
$ws.Range("V9").Value = -5986.23
$ws.Range("W9").Value = -4420.39
$ws.Range("X9").Value = -5800.75
$ws.Range("Y9").Value = -9447.82
$ws.Range("Z9").Value = -127472.18

# Row 10
$ws.Range("B10").Value = -2368.8
$ws.Range("C10").Value = -1628.42
$ws.Range("D10").Value = -4516.95
$ws.Range("E10").Value = -1574.78
$ws.Range("F10").Value = -938.27
$ws.Range("G10").Value = -2642.03
$ws.Range("H10").Value = -1076.51
$ws.Range("I10").Value = -1814.46
$ws.Range("J10").Value = -4600.93
$ws.Range("K10").Value = -443.17
$ws.Range("L10").Value = -1855.42
$ws.Range("M10").Value = -4184.69
$ws.Range("N10").Value = -2055.3
$ws.Range("O10").Value = -3698.58
$ws.Range("P10").Value = -827.59
$ws.Range("Q10").Value = -2227.42
$ws.Range("R10").Value = -3765.38
$ws.Range("S10").Value = -1564.3
$ws.Range("T10").Value = -3714.51
$ws.Range("U10").Value = -2967.24
$ws.Range("V10").Value = -3148.61
$ws.Range("W10").Value = -1744.46
$ws.Range("X10").Value = -752.22
$ws.Range("Y10").Value = -3464.47
$ws.Range("Z10").Value = -57574.51

# Row 11
$ws.Range("B11").Value = -1112.56
$ws.Range("C11").Value = -507.24
$ws.Range("D11").Value = -1874.77
$ws.Range("E11").Value = -2886.89
$ws.Range("F11").Value = -1493.62
$ws.Range("G11").Value = -673.47
$ws.Range("H11").ClearContents()
$ws.Range("I11").Value = -1331.05
$ws.Range("J11").Value = -909.25
$ws.Range("K11").Value = -507.23
$ws.Range("L11").Value = -1471.6
$ws.Range("M11").Value = -2861.88
$ws.Range("N11").Value = -509.51
$ws.Range("O11").Value = -1534.47
$ws.Range("P11").Value = -158.19
$ws.Range("Q11").Value = -3989.94
$ws.Range("R11").Value = -868.91
$ws.Range("S11").Value = -1273.64
$ws.Range("T11").Value = -1346.51
$ws.Range("U11").Value = -2356.26
$ws.Range("V11").Value = -891.23
$ws.Range("W11").Value = -3241.24
$ws.Range("X11").Value = -1142.46
$ws.Range("Y11").Value = -2595.31
$ws.Range("Z11").Value = -35537.23

# Row 12
$ws.Range("B12").Value = -22370.42
$ws.Range("C12").Value = -13198.65
$ws.Range("D12").Value = -17895.6
$ws.Range("E12").Value = -24067.3
$ws.Range("F12").Value = -17023.16
$ws.Range("G12").Value = -15068.43
$ws.Range("H12").Value = -17782.58
$ws.Range("I12").Value = -26801.76
$ws.Range("J12").Value = -17841.34
$ws.Range("K12").Value = -20835.13
$ws.Range("L12").Value = -17464.27
$ws.Range("M12").Value = -26803.79
$ws.Range("N12").Value = -19836.4
$ws.Range("O12").Value = -20582.58
$ws.Range("P12").Value = -20285.37
$ws.Range("Q12").Value = -26149.68
$ws.Range("R12").Value = -19884.91
$ws.Range("S12").Value = -22629.02
$ws.Range("T12").Value = -26543.81
$ws.Range("U12").Value = -24831.62
$ws.Range("V12").Value = -28318.63
$ws.Range("W12").Value = -16614.01
$ws.Range("X12").Value = -14368.16
$ws.Range("Y12").Value = -25558.65
$ws.Range("Z12").Value = -502755.27
